$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Dark theme completed"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = Get-Date -Year 2024 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("D9").Value = "Implemented dark theme through all app."

$ws.Range("C10").Select()
